$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "muni"
$ws.Range("K5").Select()
